$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# Update the "value" column (C2:C26) from 0.25 to the 10% quantile (0.1)
$ws.Range("C2:C26").Value = 0.1

# Update the active cell / selection to C13, matching the saved view state
$ws.Range("C13").Select()
